$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "no_polisi" (L2) cell value - the license plate number changed
$ws.Range("L2").Value = "R 5572 HS"

# Widen column B to fit the new content (closest value the engine's
# pixel-quantized ColumnWidth setter can reach to the authored 15.7109375)
$ws.Columns.Item(2).ColumnWidth = 14.84

# Reset the view: selection moves to B1 (and the sheet scrolls back to
# show column A again, clearing the previous topLeftCell="B1"/O3 selection)
$ws.Range("B1").Select()
